# Generate Report for Handoff
# Updates the "b.md" row across the Overview, zh-cn and de-de sheets to
# reflect that the b.md handoff package is ready (status + new handoff
# file names / timestamps + a version mismatch error detail).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/51b4ad276a33662cc1b6d7af155f8b0889f1c1d4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee10795f677276918208465c3bfa36dff0f7eefd/e2e/b.md."

# ---- Overview sheet: row 3 is the b.md entry ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = "2016-08-20 12:41:24"

# ---- zh-cn sheet: row 3 is the b.md entry ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-20 12:41:20"
$wsZhCn.Range("P3").Value = $errorDetail
# Excel's ColumnWidth (character units) adds ~0.83 of padding when
# serialized to the raw OOXML "width" attribute, so 39.1 round-trips to
# exactly width="40" in the saved XML.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1

# ---- de-de sheet: row 3 is the b.md entry ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-20 12:41:24"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1
